# loss function changed to MSE + center distance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

# Row 19: update the "Изменения в модели" description and fill in results / commit id
# (values are written in the same order the author filled them in, so that
# newly-created shared-string entries land in the same index order)
$ws.Range("I19").Value = "1d4cdba"
$ws.Range("H19").Value = "Train IoU: 0.53, Val IoU: 0.51. Нужно поэкспериментировать с параметрами."
$ws.Range("F19").Value = "Изменение функции потерь на комбинацию MSE и IoU, в пропорции 70 на 30"
$ws.Rows.Item(19).RowHeight = 45

# Row 20: new test entry (MSE/IoU = 20/80)
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 40
$ws.Range("D20").Value = 13
$ws.Range("H20").Value = "Train IoU: 0.47, Val IoU: 0.46. "
$ws.Range("F20").Value = "MSE/IoU = 20/80"
$ws.Range("G20").Value = "параметры теста 4"

# Row 21: new test entry (MSE/IoU = 10/90)
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 13
$ws.Range("F21").Value = "MSE/IoU = 10/90"
$ws.Range("G21").Value = "параметры теста 4"

# Update sheet view (scroll position + selection)
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("F21").Select()
